# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only the "K" column (column G) values change for this sheet's data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(1, 2, 1, 2, 2, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
